# Generate Report for Handoff
# b.md has now been handed off for translation: update its status from
# "Handed back: in sync with en-US" to "Ready for handoff" on every sheet,
# bump the handoff timestamps, and point the zh-cn / de-de handoff files at
# the newly generated xlf (hash 63290e5768f688058c7b37413b0a5c26c308f864).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the b.md row ---------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-23-12 06:23:20"

# --- zh-cn sheet: row 3 is the b.md row -------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-12 06:23:17"

# --- de-de sheet: row 3 is the b.md row -------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-12 06:23:20"
